$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "prima facie diffic|ult to understand" -> merge into a single run,
#    removing the _GoBack bookmark that sat between "diffic" and "ult".
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(" diffic" + "ult to understand.", $true, $false, $false, $false, $false, $true, 1, $false, " difficult to understand.", 2)

# ---------------------------------------------------------------------------
# 2. "...intuition behind Gibbs sampling." -> "...intuition behind sampling
#    from the posterior distribution:" and restructure the two video blocks
#    that follow (Syllabus: / Gibbs sampling: ) into hyperlinked videos.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("an explanation of the intuition behind Gibbs sampling.", $true, $false, $false, $false, $false, $true, 1, $false, "an explanation of the intuition behind sampling from the posterior distribution:", 2)

# Insert a blank paragraph right after that paragraph (before "Syllabus:").
$pFurther = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "posterior distribution:") {
        $pFurther = $para
        break
    }
}
$pFurther.Range.InsertParagraphAfter()

# Turn the "Gibbs sampling: " paragraph into a hyperlink to the first video.
$pGibbs = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "^Gibbs sampling: ") {
        $pGibbs = $para
        break
    }
}
$rGibbs = $pGibbs.Range
$rGibbs.End = $rGibbs.End - 1
$rGibbs.Text = "https://www.youtube.com/watch?v=VQRuoCawevE"
$null = $d.Hyperlinks.Add($rGibbs, "https://www.youtube.com/watch?v=VQRuoCawevE", "", "", $rGibbs.Text)

# Insert a new paragraph after it for "Bayesian posterior distribution sampling:"
$pVideo1 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "VQRuoCawevE") {
        $pVideo1 = $para
        break
    }
}
$pVideo1.Range.InsertParagraphAfter()
$pLabel2 = $pVideo1.Next()
$pLabel2.Range.InsertBefore("Bayesian posterior distribution sampling:")

# Insert another new paragraph after that for the second hyperlinked video.
$pLabel2.Range.InsertParagraphAfter()
$pVideo2 = $pLabel2.Next()
$rVideo2 = $pVideo2.Range
$rVideo2.End = $rVideo2.End - 1
$rVideo2.Text = "https://www.youtube.com/watch?v=EHqU9LE9tg8"
$null = $d.Hyperlinks.Add($rVideo2, "https://www.youtube.com/watch?v=EHqU9LE9tg8", "", "", $rVideo2.Text)

# ---------------------------------------------------------------------------
# 3. "Diagrams: ... learn visually." paragraph loses its bold paragraph
#    mark, which moves (with the trailing manual line break) into a new
#    "Code snippets:" paragraph inserted before "Market".
# ---------------------------------------------------------------------------
$pVisually = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -match "earn visually") {
        $pVisually = $para
        break
    }
}
$tVisually = $pVisually.Range.Text
$splitPos = $pVisually.Range.Start + $tVisually.Length - 2
$rSplit = $d.Range($splitPos, $splitPos)
$rSplit.InsertParagraphAfter()

$pCode = $pVisually.Next()
$insPos = $pCode.Range.Start
$codeText = "Code snippets: The reader will be provided with code snippets (which will be short) in the chapters themselves. The full code will either be provided by the companion website, or at the end of the chapters/book."
$rIns = $d.Range($insPos, $insPos)
$rIns.InsertBefore($codeText)
$rNonBold = $d.Range($insPos + 14, $insPos + $codeText.Length)
$rNonBold.Font.Bold = 0

# ---------------------------------------------------------------------------
# 4. "This title of this text seems suggestive of fulfilling a role" ->
#    "...suggestive of occupying a niche which..."
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("fulfilling a role", $true, $false, $false, $false, $false, $true, 1, $false, "occupying a niche", 2)

# ---------------------------------------------------------------------------
# 5. "6-8 months: " -> "6-10 months: "
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute("6-8 months: ", $true, $false, $false, $false, $false, $true, 1, $false, "6-10 months: ", 2)
